$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
}

# Change 1: Activation date
Replace-Text "Ativação: 01/01/2021" "Ativação: 01/01/2024"

# Change 2: Portuguese "Programa" paragraph
Replace-Text "Características Essenciais e diferenciadoras de Serviços. Ciclo de Serviços.Projeto de Novo Serviço: Planejamento estratégico (Forças de Porter, Posicionamento Estratégico). Concepção do Serviço(Conceito de Serviço. Benchmarking. SERVQUAL. Geração e Seleção de Ideias. Pacote de Serviços. Especificações deServiço). Processos (Blue Print. Padronização. Entrega do Serviço. Recrutamento e Treinamento). Instalações (Seleção eLocalização. Gestão de Evidências Físicas. Projeto do Espaço Físico. Estudo da Capacidade Produtiva). Avaliação e Melhoria" "Características Essenciais e diferenciadoras de Serviços. Ciclo de Serviços. Projeto de Novo Serviço: Planejamento estratégico (Forças de Porter, Posicionamento Estratégico). Concepção do Serviço (Conceito de Serviço. Benchmarking. SERVQUAL. Geração e Seleção de Ideias. Pacote de Serviços. Especificações de Serviço). Processos (Blue Print. Padronização. Entrega do Serviço. Recrutamento e Treinamento). Instalações (Seleção e Localização. Gestão de Evidências Físicas. Projeto do Espaço Físico. Estudo da Capacidade Produtiva). Avaliação e Melhoria. Desenvolvimento de atividade prática extensionista junto à micro e pequenos empreendedores da região (proposta de melhoria de um serviço existente). Visita (viagem didática complementar) a uma empresa prestadora de serviços para compreender o funcionamento de um sistema produtivo neste tipo de empresa."

# Change 3: English "Programa" paragraph
Replace-Text "Introduction to Service Industry. Essential and Differential Services Features. Cycle of Services. New Service Design: StrategicPlanning (Porter Forces, Strategic Positioning). Service Conception (Concept of Service, Benchmarking, SERVQUAL, Generationand Selection of Ideas, Service Package, Service Specifications). Processes (Blue Print, Standardization, Service Delivery,Recruitment and Training). Facilities (Selection and Location, Physical Evidence Management, Physical Space Project,Productive Capacity Study). Evaluation and Improvement." "Introduction to Service Industry. Essential and Differential Services Features. Cycle of Services. New Service Design: Strategic Planning (Porter Forces, Strategic Positioning). Service Conception (Concept of Service, Benchmarking, SERVQUAL, Generation and Selection of Ideas, Service Package, Service Specifications). Processes (Blue Print, Standardization, Service Delivery, Recruitment and Training). Facilities (Selection and Location, Physical Evidence Management, Physical Space Project, Productive Capacity Study). Evaluation and Improvement. Development of practical extension activities with micro and small entrepreneurs in the region (proposal to improve an existing service). Visit (complementary didactic trip) to a service provider company to understand the operation of a production system in this type of company."

# Change 4: Método
Replace-Text "Aulas Expositivas; trabalhos em grupo; exercícios individuais e palestras." "Aulas expositivas e dialogadas; dinâmicas, projetos e trabalhos em grupo; exercícios individuais; e, seminários, debates e palestras."

# Change 5: Critério
Replace-Text "Avaliação individual (Peso entre 20-40%) e do projeto realizado em equipe (peso entre 60-80%)" "Média Aritmética dos Projetos, Trabalhos, Exercícios e outras atividades avaliativas realizadas no decorrer da disciplina, considerando as questões relativas às Competências (Conhecimento, Habilidade e Atitude, que incluem a presença e participação dos alunos nas aulas) desenvolvidas"

# Change 6: Norma de recuperação
Replace-Text "NF = (MF + PR)/2, onde NF é a média final da segunda avaliação, MF é a média final da primeira avaliação e PR é a nota darecuperação." "NF = (MF + PR)/2, onde MF é a média final da avaliação e PR é uma prova de recuperação"
